$d = $word.ActiveDocument

# 1. Clear the table cell contents (A..L become empty paragraphs, matching the
#    already-empty 4th row) while keeping an (empty) paragraph per cell.
$t = $d.Tables.Item(1)
for ($r = 1; $r -le 3; $r++) {
    for ($c = 1; $c -le 4; $c++) {
        $cell = $t.Cell($r, $c)
        $full = $cell.Range
        if ($full.End - 1 -gt $full.Start) {
            $rng = $d.Range($full.Start, $full.End - 1)
            $rng.Delete()
        }
    }
}

# 2. Remove the paragraph (after the table) that used to hold the
#    bookmarkStart/bookmarkEnd pair for "_GoBack".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Bookmarks.Count -gt 0) {
        $hasGoBack = $false
        foreach ($bk in $p.Range.Bookmarks) {
            if ($bk.Name -eq "_GoBack") { $hasGoBack = $true }
        }
        if ($hasGoBack -and $p.Range.Text -eq "") {
            $p.Range.Delete()
            break
        }
    }
}

# 3. Re-add the "_GoBack" bookmark right after the "m:userdoc 'zone1'" field
#    (immediately before the end of that paragraph).
$target = $d.Paragraphs.Item(2)
$bmRange = $target.Range.Duplicate
$bmRange.Collapse(0) | Out-Null
$bmRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
